$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2000
$ws.Cells.Item(17, 10).Value = 2000
$ws.Cells.Item(17, 12).Value = 6000
$ws.Cells.Item(17, 14).Value = -6336
$ws.Cells.Item(40, 8).Value = 2374.5
$ws.Cells.Item(40, 9).Value = 2000
$ws.Cells.Item(40, 10).Value = 2642
$ws.Cells.Item(40, 11).Value = 2000
$ws.Cells.Item(40, 12).Value = 2642
$ws.Cells.Item(40, 13).Value = -1825
$ws.Cells.Item(40, 14).Value = -2992
$ws.Cells.Item(86, 8).Value = 106289.8
$ws.Cells.Item(86, 10).Value = 147271.28
$ws.Cells.Item(86, 12).Value = 147271.28
$ws.Cells.Item(86, 14).Value = -149517.28
$ws.Cells.Item(89, 8).Value = 106289.8
$ws.Cells.Item(89, 10).Value = 147271.28
$ws.Cells.Item(89, 12).Value = 736356.4
$ws.Cells.Item(89, 14).Value = -747588.4
$ws.Cells.Item(107, 8).Value = 1089.8462
$ws.Cells.Item(107, 10).Value = 450.2
$ws.Cells.Item(107, 12).Value = 450.2
$ws.Cells.Item(107, 14).Value = -4290.2
$ws.Cells.Item(129, 8).Value = 1795.7894
$ws.Cells.Item(129, 9).Value = 945.7143
$ws.Cells.Item(129, 11).Value = 2837.1429
$ws.Cells.Item(129, 13).Value = 2162.8571
$ws.Cells.Item(141, 8).Value = 4568.769
$ws.Cells.Item(141, 9).Value = 4349.5
$ws.Cells.Item(141, 10).Value = 5299.6665
$ws.Cells.Item(141, 11).Value = 13048.5
$ws.Cells.Item(141, 12).Value = 15898.9995
$ws.Cells.Item(141, 13).Value = -7868.5
$ws.Cells.Item(141, 14).Value = -26258.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(50, 8).Value = 4683
$ws.Cells.Item(50, 9).Value = 6574.5
$ws.Cells.Item(50, 10).Value = 900
$ws.Cells.Item(50, 11).Value = 6574.5
$ws.Cells.Item(50, 12).Value = 900
$ws.Cells.Item(50, 13).Value = -5860.5
$ws.Cells.Item(50, 14).Value = -2328
$ws.Cells.Item(74, 8).Value = 1000
$ws.Cells.Item(74, 9).Value = 1000
$ws.Cells.Item(74, 11).Value = 1000
$ws.Cells.Item(74, 13).Value = -126
$ws.Cells.Item(77, 8).Value = 1000
$ws.Cells.Item(77, 9).Value = 1000
$ws.Cells.Item(77, 11).Value = 5000
$ws.Cells.Item(77, 13).Value = -632
$ws.Cells.Item(110, 8).Value = 1771.05
$ws.Cells.Item(110, 9).Value = 1317.2858
$ws.Cells.Item(110, 11).Value = 1317.2858
$ws.Cells.Item(110, 13).Value = 727.7141999999999
$ws.Cells.Item(132, 8).Value = 1814.1621
$ws.Cells.Item(132, 9).Value = 1844.871
$ws.Cells.Item(132, 11).Value = 5534.613
$ws.Cells.Item(132, 13).Value = -3004.613

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1389.2
$ws.Cells.Item(94, 9).Value = 1389.2
$ws.Cells.Item(94, 11).Value = 1389.2
$ws.Cells.Item(94, 13).Value = -938.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 70673.234
$ws.Cells.Item(22, 9).Value = 131072.11
$ws.Cells.Item(22, 10).Value = 2724.5
$ws.Cells.Item(22, 11).Value = 131072.11
$ws.Cells.Item(22, 12).Value = 2724.5
$ws.Cells.Item(22, 13).Value = -130722.11
$ws.Cells.Item(22, 14).Value = -3424.5
$ws.Cells.Item(86, 8).Value = 4999
$ws.Cells.Item(86, 9).Value = 4999
$ws.Cells.Item(86, 11).Value = 4999
$ws.Cells.Item(86, 13).Value = -3876
$ws.Cells.Item(88, 8).Value = 28061.8
$ws.Cells.Item(88, 10).Value = 28499.5
$ws.Cells.Item(88, 12).Value = 28499.5
$ws.Cells.Item(88, 14).Value = -29311.5
$ws.Cells.Item(89, 8).Value = 4999
$ws.Cells.Item(89, 9).Value = 4999
$ws.Cells.Item(89, 11).Value = 24995
$ws.Cells.Item(89, 13).Value = -19379
$ws.Cells.Item(91, 8).Value = 28061.8
$ws.Cells.Item(91, 10).Value = 28499.5
$ws.Cells.Item(91, 12).Value = 28499.5
$ws.Cells.Item(91, 14).Value = -31307.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 215.45
$ws.Cells.Item(2, 9).Value = 36.42857
$ws.Cells.Item(2, 10).Value = 311.84616
$ws.Cells.Item(2, 11).Value = 218.57142
$ws.Cells.Item(2, 12).Value = 1871.07696
$ws.Cells.Item(2, 13).Value = -105.57142
$ws.Cells.Item(2, 14).Value = -2097.07696
$ws.Cells.Item(68, 8).Value = 777
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 13).Value = $null
$ws.Cells.Item(71, 8).Value = 777
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 13).Value = $null
$ws.Cells.Item(121, 8).Value = 996.6875
$ws.Cells.Item(121, 9).Value = 411.5
$ws.Cells.Item(121, 10).Value = 1191.75
$ws.Cells.Item(121, 11).Value = 1234.5
$ws.Cells.Item(121, 12).Value = 3575.25
$ws.Cells.Item(121, 13).Value = 75.5
$ws.Cells.Item(121, 14).Value = -6195.25
$ws.Cells.Item(131, 8).Value = 3235.4
$ws.Cells.Item(131, 9).Value = 1536
$ws.Cells.Item(131, 10).Value = 10033
$ws.Cells.Item(131, 11).Value = 4608
$ws.Cells.Item(131, 12).Value = 30099
$ws.Cells.Item(131, 13).Value = 432
$ws.Cells.Item(131, 14).Value = -40179
$ws.Cells.Item(134, 8).Value = 1749
$ws.Cells.Item(134, 9).Value = 1749
$ws.Cells.Item(134, 11).Value = 5247
$ws.Cells.Item(134, 13).Value = -177

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 12543.272
$ws.Cells.Item(70, 9).Value = 12194.2
$ws.Cells.Item(70, 10).Value = 12834.167
$ws.Cells.Item(70, 11).Value = 12194.2
$ws.Cells.Item(70, 12).Value = 12834.167
$ws.Cells.Item(70, 13).Value = -11924.2
$ws.Cells.Item(70, 14).Value = -13374.167
$ws.Cells.Item(73, 8).Value = 12543.272
$ws.Cells.Item(73, 9).Value = 12194.2
$ws.Cells.Item(73, 10).Value = 12834.167
$ws.Cells.Item(73, 11).Value = 12194.2
$ws.Cells.Item(73, 12).Value = 12834.167
$ws.Cells.Item(73, 13).Value = -11258.2
$ws.Cells.Item(73, 14).Value = -14706.167

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).Value = $null
$ws.Cells.Item(29, 8).Value = 59500
$ws.Cells.Item(29, 9).Value = 59000
$ws.Cells.Item(29, 10).Value = 60000
$ws.Cells.Item(29, 11).Value = 59000
$ws.Cells.Item(29, 12).Value = 60000
$ws.Cells.Item(29, 13).Value = -58705
$ws.Cells.Item(29, 14).Value = -60590
$ws.Cells.Item(46, 8).Value = 3210.5881
$ws.Cells.Item(46, 9).Value = 2398.75
$ws.Cells.Item(46, 10).Value = 3932.2222
$ws.Cells.Item(46, 11).Value = 2398.75
$ws.Cells.Item(46, 12).Value = 3932.2222
$ws.Cells.Item(46, 13).Value = -2210.75
$ws.Cells.Item(46, 14).Value = -4308.2222
$ws.Cells.Item(68, 8).Value = 4000
$ws.Cells.Item(68, 10).Value = 4000
$ws.Cells.Item(68, 12).Value = 4000
$ws.Cells.Item(68, 14).Value = -5498
$ws.Cells.Item(71, 8).Value = 4000
$ws.Cells.Item(71, 10).Value = 4000
$ws.Cells.Item(71, 12).Value = 20000
$ws.Cells.Item(71, 14).Value = -27488
$ws.Cells.Item(100, 8).Value = 3421.7
$ws.Cells.Item(100, 9).Value = 3469.4443
$ws.Cells.Item(100, 10).Value = 2992
$ws.Cells.Item(100, 11).Value = 3469.4443
$ws.Cells.Item(100, 12).Value = 2992
$ws.Cells.Item(100, 13).Value = -2928.4443
$ws.Cells.Item(100, 14).Value = -4074
$ws.Cells.Item(122, 8).Value = 5006.9644
$ws.Cells.Item(122, 9).Value = 3587.375
$ws.Cells.Item(122, 10).Value = 5574.8
$ws.Cells.Item(122, 11).Value = 10762.125
$ws.Cells.Item(122, 12).Value = 16724.4
$ws.Cells.Item(122, 13).Value = -8312.125
$ws.Cells.Item(122, 14).Value = -21624.4
$ws.Cells.Item(132, 8).Value = 102550.8
$ws.Cells.Item(132, 9).Value = 102550.8
$ws.Cells.Item(132, 11).Value = 307652.4
$ws.Cells.Item(132, 13).Value = -305122.4
$ws.Cells.Item(136, 8).Value = 6364.4287
$ws.Cells.Item(136, 9).Value = 6504
$ws.Cells.Item(136, 11).Value = 19512
$ws.Cells.Item(136, 13).Value = -16962

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(37, 8).Value = 35000
$ws.Cells.Item(37, 10).Value = 35000
$ws.Cells.Item(37, 12).Value = 35000
$ws.Cells.Item(37, 14).Value = -35406
$ws.Cells.Item(52, 8).Value = 48500
$ws.Cells.Item(52, 9).Value = 48500
$ws.Cells.Item(52, 11).Value = 48500
$ws.Cells.Item(52, 13).Value = -48274
$ws.Cells.Item(64, 8).Value = 66000
$ws.Cells.Item(64, 10).Value = 66000
$ws.Cells.Item(64, 12).Value = 66000
$ws.Cells.Item(64, 14).Value = -66496
$ws.Cells.Item(67, 8).Value = 66000
$ws.Cells.Item(67, 10).Value = 66000
$ws.Cells.Item(67, 12).Value = 66000
$ws.Cells.Item(67, 14).Value = -67716
$ws.Cells.Item(80, 8).Value = 45000
$ws.Cells.Item(80, 9).Value = 45000
$ws.Cells.Item(80, 11).Value = 45000
$ws.Cells.Item(80, 13).Value = -44002
$ws.Cells.Item(83, 8).Value = 45000
$ws.Cells.Item(83, 9).Value = 45000
$ws.Cells.Item(83, 11).Value = 135000
$ws.Cells.Item(83, 13).Value = -130008
$ws.Cells.Item(132, 8).Value = 3177.6667
$ws.Cells.Item(132, 10).Value = 3599
$ws.Cells.Item(132, 11).Value = 3599
$ws.Cells.Item(132, 12).Value = 10797
$ws.Cells.Item(132, 14).Value = -15857
